{"js": "// The document has a paragraph made of three runs that together read:\n//   \"<id>\" + \"p108r_2\" + \"</id>\"\n// The middle run (\"p108r_2\") carries plain/default formatting while the\n// outer two runs carry Courier New / color 7f6000 / size 9pt formatting.\n// The edit merges all three runs into a single run (keeping the first\n// run's formatting) whose text is the full concatenation \"<id>p108r_2</id>\".\n\nconst body = context.document.body;\n\n// Locate the run that holds the inner \"p108r_2\" text so we can get to its\n// paragraph without relying on a hard-coded paragraph index.\nconst results = body.search(\"p108r_2\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find paragraph containing 'p108r_2'\");\n}\n\nconst targetParagraph = results.items[0].paragraphs.getFirst();\n\n// Replacing the whole paragraph range's text collapses it down to a single\n// run, which adopts the formatting of the range's leading run (the\n// \"<id>\" run: Courier New / 7f6000 / 9pt) \u2014 exactly matching the merge\n// performed in the diff.\nconst paragraphRange = targetParagraph.getRange();\nparagraphRange.insertText(\"<id>p108r_2</id>\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The document contains a paragraph built from three separate runs whose\n# text concatenates to \"<id>p108r_2</id>\":\n#   run1 \"<id>\"    -> Courier New, color 7f6000, size 9pt (sz=18 half-points)\n#   run2 \"p108r_2\" -> default/plain formatting\n#   run3 \"</id>\"   -> Courier New, color 7f6000, size 9pt\n# The edit merges these three runs into a single run (carrying run1's\n# formatting) containing the full text \"<id>p108r_2</id>\".\n\n$d = $word.ActiveDocument\n\n# Find the \"p108r_2\" text anywhere in the document body and let $range\n# collapse onto the matched span (standard Word COM Find behavior).\n$range = $d.Content\n$found = $range.Find.Execute(\"p108r_2\")\nif (-not $found) {\n    throw \"Could not find 'p108r_2' in the document\"\n}\n\n# Grow the found range out to the full enclosing paragraph (wdParagraph = 4),\n# which picks up the surrounding \"<id>\" / \"</id>\" runs too.\n[void]$range.Expand(4)\n\n# Re-assigning .Text replaces the whole paragraph's run content with a\n# single new run. Word gives that new run the formatting of the first\n# character of the original range, i.e. the \"<id>\" run's Courier New /\n# 7f6000 / 9pt formatting - exactly matching the merge in the target diff.\n$range.Text = \"<id>p108r_2</id>\"\n"}
